$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 118 ("Feria Lagunitas
# de Puerto Montt" / Pomelo series). This pushes the former rows 118-203 down
# to 120-205, matching the diff's row-shift pattern, and grows the sheet
# dimension from A1:T203 to A1:T205.
$ws.Rows("118:119").Insert()

# New row 118: weekly Pomelo "Start Ruby" / "Primera" record.
$ws.Range("A118").Value = 4
$ws.Range("B118").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C118").Value = "Los Lagos"
$ws.Range("D118").Value = 44574
$ws.Range("E118").Value = 10
$ws.Range("F118").Value = "Fruta"
$ws.Range("G118").Value = 100102
$ws.Range("H118").Value = "Cítricos"
$ws.Range("I118").Value = 100102006
$ws.Range("J118").Value = "Pomelo"
$ws.Range("K118").Value = "Start Ruby"
$ws.Range("L118").Value = "Primera"
$ws.Range("M118").Value = 120
$ws.Range("N118").Value = 14000
$ws.Range("O118").Value = 15000
$ws.Range("P118").Value = 14500
$ws.Range("Q118").Value = "$/caja 14 kilos empedrada"
$ws.Range("R118").Value = "Región de O'Higgins"
$ws.Range("S118").Value = 1036
$ws.Range("T118").Value = 14

# New row 119: weekly Pomelo "Start Ruby" / "Segunda" record.
$ws.Range("A119").Value = 4
$ws.Range("B119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C119").Value = "Los Lagos"
$ws.Range("D119").Value = 44574
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100102
$ws.Range("H119").Value = "Cítricos"
$ws.Range("I119").Value = 100102006
$ws.Range("J119").Value = "Pomelo"
$ws.Range("K119").Value = "Start Ruby"
$ws.Range("L119").Value = "Segunda"
$ws.Range("M119").Value = 60
$ws.Range("N119").Value = 10000
$ws.Range("O119").Value = 10000
$ws.Range("P119").Value = 10000
$ws.Range("Q119").Value = "$/caja 14 kilos empedrada"
$ws.Range("R119").Value = "Región de O'Higgins"
$ws.Range("S119").Value = 714
$ws.Range("T119").Value = 14
